$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.641.16'
$ws.Range("E2").Value = '  -3.18%  '
$ws.Range("D3").Value = '2.296.59'
$ws.Range("E3").Value = '  -3.12%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '496.41'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '127.71'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.55%  '
$ws.Range("E7").Value = '  +0.39%  '
$ws.Range("E8").Value = '  -2.18%  '
$ws.Range("D9").Value = '2.296.44'
$ws.Range("E9").Value = '  -4.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0952'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.87%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.151'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.11%  '
$ws.Range("E12").Value = '  +0.02%  '
$ws.Range("E13").Value = '  -4.57%  '
$ws.Range("D14").Value = '2.699.75'
$ws.Range("E14").Value = '  -3.52%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.69'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.50%  '
$ws.Range("D16").Value = '54.524.12'
$ws.Range("E16").Value = '  -3.34%  '
$ws.Range("E17").Value = '  -2.63%  '
$ws.Range("D18").Value = '2.283.61'
$ws.Range("E18").Value = '  -2.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.03'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.04%  '
$ws.Range("E20").Value = '  -0.39%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '304.87'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.50'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.31%  '
$ws.Range("E23").Value = '  +0.31%  '
$ws.Range("E24").Value = '  -2.66%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '63.05'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.87%  '
$ws.Range("E26").Value = '  +0.53%  '
$ws.Range("E27").Value = '  -0.34%  '
$ws.Range("E28").Value = '  +1.38%  '
$ws.Range("D29").Value = '2.392.56'
$ws.Range("E29").Value = '  -3.75%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.12'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.45%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '170.70'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.06%  '
$ws.Range("E32").Value = '  -2.98%  '
$ws.Range("D33").Value = '0.0₃0689'
$ws.Range("E33").Value = '  -4.54%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.88'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.16%  '
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.998'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.36%  '
$ws.Range("E37").Value = '  -2.66%  '
$ws.Range("E38").Value = '  -1.21%  '
$ws.Range("E39").Value = '  +0.28%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.867'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.68%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.66'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.18%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '35.57'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.61%  '
$ws.Range("E43").Value = '  -0.46%  '
$ws.Range("E44").Value = '  -2.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '129.81'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.07%  '
$ws.Range("E46").Value = '  -1.70%  '
$ws.Range("E47").Value = '  -4.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0896'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.62%  '
$ws.Range("E49").Value = '  -2.80%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '242.79'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0482'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.12%  '
